$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (coin names, links, volume labels) -- plain assignment
$textValues = @{
    'B10' = 'One'
    'C10' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E10' = '9OneONEBestin24h'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E11' = '10WazirXWRX'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'E12' = '11MandalaExchangeTokenMDX'
    'B13' = 'LiechtensteinCryptoassetsExchange'
    'C13' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E13' = '12LiechtensteinCryptoassetsExchangeLCX'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E14' = '13BitrueCoinBTR'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E15' = '14BitMartTokenBMX'
    'B16' = 'MCDex'
    'C16' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'E16' = '15MCDexMCB'
    'B17' = 'BitForexToken'
    'C17' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E17' = '16BitForexTokenBF'
    'B18' = 'CoinExToken'
    'C18' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E18' = '17CoinExTokenCET'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'E41' = '40KickTokenKICK'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E42' = '41BKEXTokenBKK'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'E43' = '42CEJICEJI'
    'E47' = '46CoinbaseStockTokenCOIN'
    'E48' = '47BOLOBOLOWorstin24h'
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

# Numeric-looking price strings -- force Text format so values like
# "245.20" / "0.001590" keep their exact digits/trailing zeros instead
# of being parsed into numbers.
$priceValues = @{
    'D2' = '245.20'
    'D4' = '5.386'
    'D5' = '0.05978'
    'D6' = '3.395'
    'D7' = '6.401'
    'D8' = '0.8110'
    'D9' = '0.9598'
    'D10' = '0.01123'
    'D11' = '0.1428'
    'D12' = '0.07415'
    'D13' = '0.03434'
    'D14' = '0.03065'
    'D15' = '0.09417'
    'D16' = '4.002'
    'D17' = '0.001590'
    'D18' = '0.04818'
    'D19' = '0.006146'
    'D21' = '0.0009884'
    'D23' = '3.744'
    'D24' = '2.187'
    'D27' = '0.0002463'
    'D40' = '0.03960'
    'D41' = '0.006504'
    'D42' = '0.1073'
    'D43' = '0.002301'
    'D44' = '0.005314'
    'D45' = '0.00005245'
    'D47' = '0.6704'
    'D48' = '0.02867'
    'D50' = '0.01011'
}

foreach ($addr in $priceValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceValues[$addr]
    $cell.NumberFormat = "General"
}
